$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95, column B currently holds text "3" -- convert it to a real number.
$ws.Range("B95").Value = 3

# Insert new annotation row 96 with the new data.
$ws.Range("A96").Value = "Ying Tang"
$ws.Range("B96").NumberFormat = "@"
$ws.Range("B96").Value = "3"
$ws.Range("C96").Value = "To address"
$ws.Range("D96").Value = "ACK"
$ws.Range("E96").Value = "EXP"
$ws.Range("F96").Value = "3419a239-823d-4d38-8055-389a9317394a"
$ws.Range("G96").Value = "SJa9iHgAZ_annotated.xlsx"
$ws.Range("H96").Value = "To address Reviewer 2 comment on iterative inference in shared Resnet, we added two sections in Appendix reporting metrics (cosine loss, accuracy, l1 ratio) on shared Resnet, and on the unrolled to more steps Resnet."
